$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Analysis notes" column (F) with header and notes for several rows.
# F85 was filled in first (matches shared-string insertion order), followed
# by the header and the rest of the notes in row order.
$ws.Range("F85").Value = "Plotted. Steps take back and forth (small reduction of steps)"
$ws.Range("F1").Value = "Analysis notes"
$ws.Range("F3").Value = "No motion"
$ws.Range("F4").Value = "Good videos, problem tracking F1F7_B9_B_4ms_40V"
$ws.Range("F5").Value = "video out of focus, no motion. 16 ms/35 V shows some motion. 32 ms /35 v back and forth motion"
$ws.Range("F6").Value = "Not very useful"
$ws.Range("F7").Value = "Good at 35 and 40 V "
$ws.Range("F8").Value = "40 V nice clear steps"
$ws.Range("F9").Value = "Useless"
$ws.Range("F10").Value = "Clear videos for some voltages. See graphs"

# Size the new column to match the widened "Analysis notes" column
$ws.Columns.Item(6).ColumnWidth = 20.2

# Update the selection / scroll position to match the latest view state
[void]$ws.Range("E19").Select()
